# Fixed naive component forecaster bug - Presentation state 11.02.
# The error-table rows shift down by one lag (row N -> row N+1 for N=2..10);
# a brand-new set of error-metric values is inserted at row 2 (the new
# first QoQ lag), and the previous last row (old row 11) is dropped.
# Column A (the lag/quarter label) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.04915694638892727
$ws.Range("C2").Value = 0.6586326132859144
$ws.Range("D2").Value = 1.297461883598808
$ws.Range("E2").Value = 1.13906184362343
$ws.Range("F2").Value = 1.169184523175352
$ws.Range("G2").Value = 19
$ws.Range("B3").Value = -0.07072591166527165
$ws.Range("C3").Value = 0.5579505877006946
$ws.Range("D3").Value = 0.6529064599314954
$ws.Range("E3").Value = 0.8080262742828945
$ws.Range("F3").Value = 0.8282610236588537
$ws.Range("G3").Value = 18
$ws.Range("B4").Value = 0.03916263344646955
$ws.Range("C4").Value = 0.552040872335158
$ws.Range("D4").Value = 0.6357226611696118
$ws.Range("E4").Value = 0.7973221815361792
$ws.Range("F4").Value = 0.8208689046632819
$ws.Range("G4").Value = 17
$ws.Range("B5").Value = 0.05955592157360488
$ws.Range("C5").Value = 0.532271854806013
$ws.Range("D5").Value = 0.7198117377579547
$ws.Range("E5").Value = 0.8484171955812509
$ws.Range("F5").Value = 0.8740799841134135
$ws.Range("G5").Value = 16
$ws.Range("B6").Value = 0.09235442246305205
$ws.Range("C6").Value = 0.5824508295849771
$ws.Range("D6").Value = 0.7506227629998017
$ws.Range("E6").Value = 0.8663848815623468
$ws.Range("F6").Value = 0.8916838547205332
$ws.Range("G6").Value = 15
$ws.Range("B7").Value = 0.1521483464337763
$ws.Range("C7").Value = 0.6635260061488867
$ws.Range("D7").Value = 1.004366658952618
$ws.Range("E7").Value = 1.002180951202236
$ws.Range("F7").Value = 1.0279571060649
$ws.Range("G7").Value = 14
$ws.Range("B8").Value = 0.159095420907867
$ws.Range("C8").Value = 0.6746913222580827
$ws.Range("D8").Value = 0.8632158529041898
$ws.Range("E8").Value = 0.9290941033631576
$ws.Range("F8").Value = 0.9527485895797485
$ws.Range("G8").Value = 13
$ws.Range("B9").Value = 0.2207410413997515
$ws.Range("C9").Value = 0.686169913841685
$ws.Range("D9").Value = 0.7776458313225046
$ws.Range("E9").Value = 0.8818422939066285
$ws.Range("F9").Value = 0.8917312420012032
$ws.Range("G9").Value = 12
$ws.Range("B10").Value = 0.220829446489595
$ws.Range("C10").Value = 0.7371449325608843
$ws.Range("D10").Value = 1.012670326230143
$ws.Range("E10").Value = 1.006315222099985
$ws.Range("F10").Value = 1.029706341620059
$ws.Range("G10").Value = 11
$ws.Range("B11").Value = 0.1936765229212045
$ws.Range("C11").Value = 0.6934653330366425
$ws.Range("D11").Value = 1.055080256735366
$ws.Range("E11").Value = 1.027170996833227
$ws.Range("F11").Value = 1.063312257473744
$ws.Range("G11").Value = 10
